$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formula")
$ws.Activate()

# Insert a new column before the old "G" (the merged "Shared formula" header
# block G4:I4, and the ROW()-based formulas in G5:I5, all shift one column
# right to H4:J4 / H5:J5). Excel's column insert copies the left neighbour's
# formatting into the new column automatically.
$ws.Columns("G").Insert()

# Give the new column the same width as the other formula columns (D:F).
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# G4 is the new header cell above the new formula column: empty but with a
# medium left/top border (no right/bottom) to visually continue the box to
# its right, inheriting the centered/wrapped alignment + "applied fill" flag
# copied in from column F's header style.
$ws.Range("G4").Borders.Item(10).LineStyle = 0   # xlEdgeRight -> none

# G5: new formula mirroring F5's D$5+E$5 pattern but without the $ row lock.
$ws.Range("G5").Formula = "=D5+E5"

# The selection that used to sit on J1 now sits one column to the left (I1)
# since the new column pushed everything over by one.
$ws.Range("I1").Select()

$wb.Save()
